# Auto-generated script applying the cryptos.xlsx diff
# (updated crypto price/volume figures, and two pairs of rows whose
# Coin/Link/Price/Volume data were swapped: rows 22<->23 and 50<->51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use Formula with a leading apostrophe (classic 'force text' prefix) so
# that numeric-looking strings (e.g. '6.61', '1.00', '0.100') are stored
# as literal text rather than being coerced into numbers by Excel, and
# without altering each cell's NumberFormat/style.

$ws.Range("D2").Formula = "'68.262.48"
$ws.Range("E2").Formula = "'  -0.45%  "
$ws.Range("D3").Formula = "'3.274.77"
$ws.Range("E3").Formula = "'  +0.82%  "
$ws.Range("E4").Formula = "'  -0.03%  "
$ws.Range("D5").Formula = "'584.42"
$ws.Range("E5").Formula = "'  +0.06%  "
$ws.Range("D6").Formula = "'185.27"
$ws.Range("E6").Formula = "'  +2.40%  "
$ws.Range("E7").Formula = "'  +0.01%  "
$ws.Range("E8").Formula = "'  +1.53%  "
$ws.Range("E9").Formula = "'  -2.68%  "
$ws.Range("D10").Formula = "'6.61"
$ws.Range("E10").Formula = "'  -0.60%  "
$ws.Range("D11").Formula = "'0.411"
$ws.Range("E11").Formula = "'  -2.35%  "
$ws.Range("D12").Formula = "'3.840.35"
$ws.Range("E12").Formula = "'  +0.72%  "
$ws.Range("E13").Formula = "'  +0.89%  "
$ws.Range("D14").Formula = "'27.61"
$ws.Range("E14").Formula = "'  -2.11%  "
$ws.Range("D15").Formula = "'68.231.69"
$ws.Range("E15").Formula = "'  -0.45%  "
$ws.Range("D16").Formula = "'0.0000169"
$ws.Range("E16").Formula = "'  -1.32%  "
$ws.Range("D17").Formula = "'3.317.92"
$ws.Range("E17").Formula = "'  +1.88%  "
$ws.Range("D18").Formula = "'5.74"
$ws.Range("E18").Formula = "'  -1.29%  "
$ws.Range("E19").Formula = "'  -0.73%  "
$ws.Range("D20").Formula = "'417.30"
$ws.Range("E20").Formula = "'  +5.98%  "
$ws.Range("D21").Formula = "'7.58"
$ws.Range("E21").Formula = "'  -0.98%  "
$ws.Range("B22").Formula = "'Litecoin"
$ws.Range("C22").Formula = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Formula = "'71.55"
$ws.Range("E22").Formula = "'  +0.27%  "
$ws.Range("B23").Formula = "'Dai"
$ws.Range("C23").Formula = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Formula = "'1.00"
$ws.Range("E23").Formula = "'  -0.09%  "
$ws.Range("D24").Formula = "'0.509"
$ws.Range("E24").Formula = "'  -1.21%  "
$ws.Range("E25").Formula = "'  -0.50%  "
$ws.Range("E26").Formula = "'  -1.27%  "
$ws.Range("D27").Formula = "'9.46"
$ws.Range("E27").Formula = "'  -1.46%  "
$ws.Range("E28").Formula = "'  +0.37%  "
$ws.Range("D29").Formula = "'1.96"
$ws.Range("E29").Formula = "'  -1.37%  "
$ws.Range("D30").Formula = "'22.76"
$ws.Range("E30").Formula = "'  -0.85%  "
$ws.Range("D31").Formula = "'5.49"
$ws.Range("E31").Formula = "'  -3.24%  "
$ws.Range("D32").Formula = "'6.90"
$ws.Range("E32").Formula = "'  -2.78%  "
$ws.Range("E33").Formula = "'  +0.04%  "
$ws.Range("D34").Formula = "'1.25"
$ws.Range("E34").Formula = "'  -1.91%  "
$ws.Range("D35").Formula = "'164.15"
$ws.Range("E35").Formula = "'  -0.10%  "
$ws.Range("E36").Formula = "'  -2.15%  "
$ws.Range("D37").Formula = "'1.90"
$ws.Range("E37").Formula = "'  -1.84%  "
$ws.Range("D38").Formula = "'27.18"
$ws.Range("E38").Formula = "'  +3.53%  "
$ws.Range("E39").Formula = "'  -2.57%  "
$ws.Range("D40").Formula = "'4.48"
$ws.Range("E40").Formula = "'  -2.32%  "
$ws.Range("D41").Formula = "'6.35"
$ws.Range("E41").Formula = "'  -3.13%  "
$ws.Range("D42").Formula = "'2.664.07"
$ws.Range("E42").Formula = "'  +2.85%  "
$ws.Range("D43").Formula = "'40.81"
$ws.Range("E43").Formula = "'  -1.46%  "
$ws.Range("E44").Formula = "'  -1.08%  "
$ws.Range("E45").Formula = "'  -1.32%  "
$ws.Range("D46").Formula = "'338.18"
$ws.Range("E46").Formula = "'  -1.46%  "
$ws.Range("D47").Formula = "'24.55"
$ws.Range("E47").Formula = "'  -0.40%  "
$ws.Range("E48").Formula = "'  -2.29%  "
$ws.Range("E49").Formula = "'  +0.00%  "
$ws.Range("B50").Formula = "'ONDO"
$ws.Range("C50").Formula = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Formula = "'0.981"
$ws.Range("E50").Formula = "'  +0.42%  "
$ws.Range("B51").Formula = "'Stellar"
$ws.Range("C51").Formula = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Formula = "'0.100"
$ws.Range("E51").Formula = "'  -1.02%  "
